$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1247905694626398
$ws.Range("C3").Value = 0.02661499342593231
$ws.Range("C4").Value = 0.0933881345806391
$ws.Range("C5").Value = 0.06643477924217901
$ws.Range("C6").Value = 0.05695587007455943
$ws.Range("C7").Value = 0.0339732674219591
$ws.Range("C8").Value = 0.1592739769840469
$ws.Range("C9").Value = 0.07245769372894408
$ws.Range("C10").Value = 0.07920502053931509
$ws.Range("C11").Value = 0.09925769745512135
$ws.Range("C12").Value = 0.08328591025254585
$ws.Range("C13").Value = 0.05339220446728885
$ws.Range("C14").Value = 0.05096988236482914
